$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows starting at row 34 (this pushes the existing
# rows 34-38 -- the two "Castle Brite" rows and "Modesto" rows -- down to
# rows 37-41, unchanged) so that new data can be inserted in rows 34-36.
$ws.Range("A34:A36").EntireRow.Insert()

# Row 34: Comercializadora del Agro de Limarí - Damasco - Dina - Especial (Illapel)
$ws.Range("A34").Value = 2
$ws.Range("B34").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C34").Value = "Coquimbo"
$ws.Range("D34").Value = 44931
$ws.Range("D34").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E34").Value = 4
$ws.Range("F34").Value = "Fruta"
$ws.Range("G34").Value = 100103
$ws.Range("H34").Value = "Frutos de hueso (carozo)"
$ws.Range("I34").Value = 100103003
$ws.Range("J34").Value = "Damasco"
$ws.Range("K34").Value = "Dina"
$ws.Range("L34").Value = "Especial"
$ws.Range("M34").Value = 300
$ws.Range("N34").Value = 22000
$ws.Range("O34").Value = 23000
$ws.Range("P34").Value = 22500
$ws.Range("Q34").Value = "$/caja 16 kilos"
$ws.Range("R34").Value = "Illapel"
$ws.Range("S34").Value = 1406
$ws.Range("T34").Value = 16

# Row 35: Comercializadora del Agro de Limarí - Damasco - Dina - Primera (Illapel)
$ws.Range("A35").Value = 2
$ws.Range("B35").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C35").Value = "Coquimbo"
$ws.Range("D35").Value = 44931
$ws.Range("D35").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E35").Value = 4
$ws.Range("F35").Value = "Fruta"
$ws.Range("G35").Value = 100103
$ws.Range("H35").Value = "Frutos de hueso (carozo)"
$ws.Range("I35").Value = 100103003
$ws.Range("J35").Value = "Damasco"
$ws.Range("K35").Value = "Dina"
$ws.Range("L35").Value = "Primera"
$ws.Range("M35").Value = 400
$ws.Range("N35").Value = 19000
$ws.Range("O35").Value = 20000
$ws.Range("P35").Value = 19500
$ws.Range("Q35").Value = "$/caja 16 kilos"
$ws.Range("R35").Value = "Illapel"
$ws.Range("S35").Value = 1219
$ws.Range("T35").Value = 16

# Row 36: Comercializadora del Agro de Limarí - Damasco - Dina - Segunda (Illapel)
$ws.Range("A36").Value = 2
$ws.Range("B36").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C36").Value = "Coquimbo"
$ws.Range("D36").Value = 44931
$ws.Range("D36").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E36").Value = 4
$ws.Range("F36").Value = "Fruta"
$ws.Range("G36").Value = 100103
$ws.Range("H36").Value = "Frutos de hueso (carozo)"
$ws.Range("I36").Value = 100103003
$ws.Range("J36").Value = "Damasco"
$ws.Range("K36").Value = "Dina"
$ws.Range("L36").Value = "Segunda"
$ws.Range("M36").Value = 400
$ws.Range("N36").Value = 15000
$ws.Range("O36").Value = 16000
$ws.Range("P36").Value = 15500
$ws.Range("Q36").Value = "$/caja 16 kilos"
$ws.Range("R36").Value = "Illapel"
$ws.Range("S36").Value = 969
$ws.Range("T36").Value = 16
